$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("FM2").Value = 74.87
$ws.Range("FP2").Value = 774.87
$ws.Range("FM3").Value = 531.89
$ws.Range("FP3").Value = 1306.76
$ws.Range("FM4").Value = -132.79
$ws.Range("FP4").Value = 1173.97
$ws.Range("FM5").Value = 302.68
$ws.Range("FP5").Value = 1476.65
$ws.Range("FM6").Value = 126.91
$ws.Range("FP6").Value = 1603.56
$ws.Range("FM7").Value = 36.46
$ws.Range("FP7").Value = 1640.02
$ws.Range("FM8").Value = 16.82
$ws.Range("FP8").Value = 1656.84
$ws.Range("FM9").Value = -28.15
$ws.Range("FP9").Value = 1628.69
$ws.Range("FM10").Value = 83.72
$ws.Range("FP10").Value = 1712.41
$ws.Range("FM11").Value = 542.87
$ws.Range("FP11").Value = 2255.28
$ws.Range("FM12").Value = -525.3200000000001
$ws.Range("FP12").Value = 1729.96
$ws.Range("FM13").Value = -338.17
$ws.Range("FP13").Value = 1391.79
$ws.Range("FM14").Value = -166.96
$ws.Range("FP14").Value = 1224.83
$ws.Range("FM15").Value = 289.05
$ws.Range("FP15").Value = 1513.88
$ws.Range("FM16").Value = 16.46
$ws.Range("FP16").Value = 1530.34
$ws.Range("FM17").Value = -58.78
$ws.Range("FP17").Value = 1471.56
$ws.Range("FM18").Value = -67.45999999999999
$ws.Range("FP18").Value = 1404.1
$ws.Range("FM19").Value = -93.84
$ws.Range("FP19").Value = 1310.26
$ws.Range("FM20").Value = 183.2
$ws.Range("FP20").Value = 1493.46
$ws.Range("FM21").Value = 3.54
$ws.Range("FP21").Value = 1497
$ws.Range("FM22").Value = -45.15
$ws.Range("FP22").Value = 1451.85
$ws.Range("FM23").Value = 531.2
$ws.Range("FP23").Value = 1983.05
$ws.Range("FM24").Value = 283.21
$ws.Range("FP24").Value = 2266.26
$ws.Range("FM25").Value = 92.04000000000001
$ws.Range("FP25").Value = 2358.3
$ws.Range("FM26").Value = 8.67
$ws.Range("FP26").Value = 2366.97
$ws.Range("FM27").Value = -57.19
$ws.Range("FP27").Value = 2309.78
$ws.Range("FM28").Value = -36.12
$ws.Range("FP28").Value = 2273.66
$ws.Range("FM29").Value = -508.87
$ws.Range("FP29").Value = 1764.79
$ws.Range("FM30").Value = -531.54
$ws.Range("FP30").Value = 1233.25
$ws.Range("FM31").Value = 109.39
$ws.Range("FP31").Value = 1342.64
$ws.Range("FM32").Value = -110.31
$ws.Range("FP32").Value = 1232.33
$ws.Range("FM33").Value = -211.41
$ws.Range("FP33").Value = 1020.92
$ws.Range("FM34").Value = -289.14
$ws.Range("FP34").Value = 731.7800000000002
$ws.Range("FM35").Value = 34.69
$ws.Range("FP35").Value = 766.4700000000003
$ws.Range("FM36").Value = -144.48
$ws.Range("FP36").Value = 621.9900000000002
$ws.Range("FM37").Value = -228.05
$ws.Range("FP37").Value = 393.9400000000002
$ws.Range("FM38").Value = 95.23
$ws.Range("FP38").Value = 489.1700000000002
$ws.Range("FM39").Value = 58.41
$ws.Range("FP39").Value = 547.5800000000003
$ws.Range("FM40").Value = -317.29
$ws.Range("FP40").Value = 230.2900000000002
$ws.Range("FM41").Value = -216.37
$ws.Range("FP41").Value = 13.92000000000024
$ws.Range("FM42").Value = -135.63
$ws.Range("FP42").Value = -121.7099999999998
$ws.Range("FM43").Value = -148.38
$ws.Range("FP43").Value = -270.0899999999997
$ws.Range("FM44").Value = 863.47
$ws.Range("FP44").Value = 593.3800000000003
$ws.Range("FM45").Value = -501.44
$ws.Range("FP45").Value = 91.94000000000034
$ws.Range("FM46").Value = -208.76
$ws.Range("FP46").Value = -116.8199999999997
$ws.Range("FM47").Value = 82.67
$ws.Range("FP47").Value = -34.14999999999965
$ws.Range("FM48").Value = -110.84
$ws.Range("FP48").Value = -144.9899999999997
$ws.Range("FM49").Value = -198.84
$ws.Range("FP49").Value = -343.8299999999997
$ws.Range("FM50").Value = -18.41
$ws.Range("FP50").Value = -362.2399999999997
$ws.Range("FM51").Value = 31.51
$ws.Range("FP51").Value = -330.7299999999997
$ws.Range("FP52").Value = -335.8599999999997
$ws.Range("FM53").Value = -226.63
$ws.Range("FP53").Value = -562.4899999999998
$ws.Range("FM54").Value = -247.88
$ws.Range("FP54").Value = -810.3699999999998
$ws.Range("FM55").Value = 5.49
$ws.Range("FP55").Value = -804.8799999999998
$ws.Range("FM56").Value = -352.16
$ws.Range("FP56").Value = -1157.04
$ws.Range("FM57").Value = -210.69
$ws.Range("FP57").Value = -1367.73
$ws.Range("FM58").Value = -89.06
$ws.Range("FP58").Value = -1456.79
$ws.Range("FM59").Value = 129.39
$ws.Range("FP59").Value = -1327.4
$ws.Range("FM60").Value = 38.94
$ws.Range("FP60").Value = -1288.46
$ws.Range("FM61").Value = -78.79000000000001
$ws.Range("FP61").Value = -1367.25
$ws.Range("FM62").Value = 47.97
$ws.Range("FP62").Value = -1319.28
$ws.Range("FM63").Value = -23.9
$ws.Range("FP63").Value = -1343.18
$ws.Range("FM64").Value = -235.84
$ws.Range("FP64").Value = -1579.02
$ws.Range("FM65").Value = -520.55
$ws.Range("FP65").Value = -2099.57
